$wb = $excel.ActiveWorkbook

# --- Update parameters sheet values (stepsize_rule + theta_set) ---
$ws = $wb.Worksheets.Item("parameters")
$ws.Range("B7").Value = "Declining"
$ws.Range("B8").Value = "1 5 10 20 50"

# --- Bump the zoom level on the Instructions sheet ---
$wsInstr = $wb.Worksheets.Item("Instructions")
$wsInstr.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 120

# --- Re-activate "parameters" as the selected tab, with B9 selected ---
$ws.Activate() | Out-Null
$ws.Range("B9").Select() | Out-Null
